$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap contents of column B and column C (rows 1-4)
$colB = $ws.Range("B1:B4").Value()
$colC = $ws.Range("C1:C4").Value()
$ws.Range("B1:B4").Value = $colC
$ws.Range("C1:C4").Value = $colB

# Swap contents of column D and column E (rows 1-4)
$colD = $ws.Range("D1:D4").Value()
$colE = $ws.Range("E1:E4").Value()
$ws.Range("D1:D4").Value = $colE
$ws.Range("E1:E4").Value = $colD

# Swap the column widths accordingly (B<->C, D<->E)
$wB = $ws.Columns.Item(2).ColumnWidth()
$wC = $ws.Columns.Item(3).ColumnWidth()
$wD = $ws.Columns.Item(4).ColumnWidth()
$wE = $ws.Columns.Item(5).ColumnWidth()
$ws.Columns.Item(2).ColumnWidth = $wC
$ws.Columns.Item(3).ColumnWidth = $wB
$ws.Columns.Item(4).ColumnWidth = $wE
$ws.Columns.Item(5).ColumnWidth = $wD

# Update the active selection to C11
[void]$ws.Range("C11").Select()
